$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.57
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.75

# Row 4
$ws.Range("I4").Value = 2.05

# Row 5
$ws.Range("G5").Value = 1.55

# Row 6
$ws.Range("I6").Value = 1.75
$ws.Range("M6").Value = 1.1
$ws.Range("N6").Value = 7

# Row 7
$ws.Range("R7").Value = 1.47

# Row 9
$ws.Range("I9").Value = 3.6
$ws.Range("K9").Value = 1.92

# Row 10
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 8

# Row 12
$ws.Range("G12").Value = 1.67
$ws.Range("H12").Value = 3.2
$ws.Range("I12").Value = 5.75
$ws.Range("J12").Value = 2.4
$ws.Range("L12").Value = 6
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 7
$ws.Range("U12").Value = 2.25
$ws.Range("V12").Value = 1.57
$ws.Range("W12").Value = 5
$ws.Range("X12").Value = 6.5
$ws.Range("Y12").Value = 9.5
$ws.Range("Z12").Value = 12
$ws.Range("AI12").Value = 26
$ws.Range("AJ12").Value = 21
$ws.Range("AK12").Value = 67
$ws.Range("AL12").Value = 51
$ws.Range("AN12").Value = 3.5
$ws.Range("AO12").Value = 9.5
$ws.Range("AT12").Value = 2.37
$ws.Range("AU12").Value = 10
$ws.Range("AW12").Value = 7
$ws.Range("AX12").Value = 34
$ws.Range("AZ12").Value = 126
$ws.Range("BA12").Value = 201

# Row 13
$ws.Range("Q13").Value = 1.87
$ws.Range("R13").Value = 1.87
